# Edit script: transform the expenses workbook per commit:
# "misc changes to expenses / added invoice and formatted spreadsheet"
#
# Original layout: row1 = category headers (lock, NES controllers, digikey,
# goodwill, ebay) across B:F, row2 = Carlos Mariscal's amounts for each
# category.
#
# New layout: row2 = person headers (Carlos, Jeff, Luis, Chelsea, Saida)
# across B:F, rows 3-10 = expense categories (one per row) with each
# person's contribution in their column. Rows 11-14 blank spacer rows,
# row15 = totals (SUM formulas), row16 blank, rows 17-20 blank (accounting
# format placeholders).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Start from a clean slate for the cell grid (keep column defs).
$ws.Cells.Clear()

$acct = "_(""$""* #,##0.00_);_(""$""* \(#,##0.00\);_(""$""* ""-""??_);_(@_)"

# ---- Header row (row 2): people ----
$ws.Range("B2").Value = "Carlos Mariscal"
$ws.Range("C2").Value = "Jeff Alcoke "
$ws.Range("D2").Value = "Luis Santiago"
$ws.Range("E2").Value = "Chelsea Throop"
$ws.Range("F2").Value = "Saida Ahker"

# G2:Z2 filler cells carry over the old "$ currency" style (style index 1
# in the original workbook, numFmtId 164).
$ws.Range("G2:Z2").NumberFormat = """$""#,##0.00"

# ---- Category rows 3-7 (existing categories, transposed) ----
$ws.Range("A3").Value = "lock"
$ws.Range("A4").Value = "NES controllers"
$ws.Range("A5").Value = "digikey"
$ws.Range("A6").Value = "goodwill (12V bricks)"
$ws.Range("A7").Value = "ebay (controllers)"

# G3:Z5 filler cells (only rows 3-5 had them originally; carried along
# with the same style as G2:Z2).
$ws.Range("G3:Z5").NumberFormat = """$""#,##0.00"

# ---- New invoice rows 8-10 (Luis Santiago's LCD-related purchases) ----
$ws.Range("A8").Value = "LockBox boards"
$ws.Range("A9").Value = "LCD Board"

# ---- Apply the accounting number format to the whole data block ----
$ws.Range("B3:F14").NumberFormat = $acct
$ws.Range("B16:F16").NumberFormat = $acct
$ws.Range("B17:F20").NumberFormat = $acct

# ---- Data values ----
$ws.Range("B3").Value = 15
$ws.Range("E3").Value = 0
$ws.Range("B4").Value = 6
$ws.Range("B5").Value = 33.83
$ws.Range("D5").Value = 17.13
$ws.Range("B6").Value = 4
$ws.Range("B7").Value = 21
$ws.Range("D8").Value = 38.3
$ws.Range("D9").Value = 26.6
$ws.Range("D10").Value = 24

# ---- Totals row 15 ----
$ws.Range("A15").Value = "total"

# ---- Last invoice row (digikey LCD) ----
$ws.Range("A10").Value = "digikey LCD"

$ws.Range("A15").Borders.Item(8).LineStyle = 1
$ws.Range("B15:F15").Borders.Item(8).LineStyle = 1
$ws.Range("B15:F15").NumberFormat = $acct
$ws.Range("B15").Formula = "=SUM(B3:B14)"
$ws.Range("C15:F15").Formula = "=SUM(C3:C14)"

# ---- Column widths / layout tweaks ----
$ws.Columns.Item(2).ColumnWidth = 15.625

Write-Host "Edit applied"
